$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (AD1:AF1) for Wins / Losses / Ties.
# Copy the formatting of the existing header style (AC1, style index "1":
# bold font, thin border, center/top alignment) onto the new header cells
# before setting their text, so they match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (82 wins, 79 losses, 0 ties) for every player
# row (2-50) in the new AD/AE/AF columns.
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 82
    $ws.Cells.Item($r, 31).Value = 79
    $ws.Cells.Item($r, 32).Value = 0
}
